$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates - Excel keeps these as text naturally
# (percentages with spaces, and D-values that contain multiple dots / non-numeric chars)
$textUpdates = @{
    "D2"  = "26.001.47"
    "E2"  = "  +0.33%  "
    "D3"  = "1.642.31"
    "E3"  = "  +0.38%  "
    "E4"  = "  +0.40%  "
    "E6"  = "  -0.08%  "
    "E7"  = "  +0.38%  "
    "E8"  = "  +0.58%  "
    "E9"  = "  +0.68%  "
    "E10" = "  -0.22%  "
    "E11" = "  +0.17%  "
    "D12" = "1.870.24"
    "E12" = "  +0.41%  "
    "E13" = "  +0.64%  "
    "D14" = "1.653.19"
    "E15" = "  +0.37%  "
    "D16" = "0.0$([char]0x2083)0764"
    "E16" = "  +1.18%  "
    "E17" = "  +1.42%  "
    "D18" = "26.027.10"
    "E18" = "  +0.37%  "
    "E19" = "  +0.40%  "
    "E20" = "  +0.37%  "
    "E21" = "  -0.77%  "
    "E22" = "  +0.08%  "
    "E23" = "  -1.03%  "
    "E24" = "  +4.47%  "
    "E25" = "  -0.93%  "
    "E26" = "  +0.41%  "
    "E27" = "  -0.50%  "
    "E28" = "  +0.60%  "
    "E30" = "  +0.55%  "
    "E31" = "  -1.10%  "
    "E32" = "  +0.21%  "
    "E33" = "  +1.74%  "
    "E34" = "  -0.55%  "
    "D37" = "1.131.16"
    "E37" = "  -0.69%  "
    "E38" = "  -0.78%  "
    "E39" = "  -0.37%  "
    "E40" = "  +0.16%  "
    "E41" = "  +0.76%  "
    "E42" = "  -0.40%  "
    "E43" = "  +0.08%  "
    "E44" = "  +3.43%  "
    "E45" = "  +0.30%  "
    "E46" = "  +3.77%  "
    "E47" = "  -1.39%  "
    "E48" = "  +1.89%  "
    "E49" = "  -0.17%  "
    "E50" = "  +0.30%  "
    "E51" = "  -0.61%  "
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Numeric-looking text updates - would be auto-converted to numbers by Excel,
# so force text format, assign, then restore default style (no custom format
# left behind on the cell).
$numericLookingUpdates = @{
    "D5"  = "215.80"
    "D7"  = "1.01"
    "D10" = "19.55"
    "D11" = "0.0795"
    "D17" = "63.44"
    "D20" = "194.21"
    "D21" = "4.36"
    "D23" = "6.21"
    "D26" = "1.01"
    "D27" = "143.01"
    "D30" = "1.25"
    "D35" = "2.47"
    "D41" = "5.47"
    "D42" = "99.08"
    "D43" = "0.798"
    "D45" = "56.53"
    "D48" = "7.77"
}

foreach ($addr in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$addr]
    $cell.Style = "Normal"
}
